{"js": "// Update the date heading paragraph (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text.trim() === \"2025-09-28 Sunday\") {\n  // Replace just the text of the run(s) in this paragraph while keeping\n  // the paragraph's own formatting (font/size) intact.\n  dateParagraph.getRange().insertText(\"2025-09-29 Monday\", Word.InsertLocation.replace);\n}\n\n// Update the 20x5 table of addition/subtraction problems. Writing the\n// `values` 2D array back preserves each cell's existing run formatting\n// (font, size, paragraph alignment) and only swaps the text content.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst newValues = [\n  [\"69-43=26\", \"33+53=86\", \"17+53=70\", \"60+36=96\", \"18+32=50\"],\n  [\"36+63=99\", \"64+29=93\", \"37-15=22\", \"97-32=65\", \"98-18=80\"],\n  [\"26+3=29\", \"37-30=7\", \"5+88=93\", \"59-15=44\", \"26-15=11\"],\n  [\"97+0=97\", \"3+93=96\", \"6+31=37\", \"77-49=28\", \"11+19=30\"],\n  [\"15+18=33\", \"2-0=2\", \"53+33=86\", \"8+20=28\", \"76+5=81\"],\n  [\"50+26=76\", \"57+2=59\", \"46+5=51\", \"33+46=79\", \"0+8=8\"],\n  [\"94-82=12\", \"42-3=39\", \"58+19=77\", \"37-13=24\", \"19+80=99\"],\n  [\"32-18=14\", \"46+5=51\", \"81-40=41\", \"30+67=97\", \"99-55=44\"],\n  [\"69+21=90\", \"56+37=93\", \"49+22=71\", \"96-5=91\", \"37-8=29\"],\n  [\"12+69=81\", \"6+2=8\", \"53+18=71\", \"35+16=51\", \"26+66=92\"],\n  [\"59-43=16\", \"85-32=53\", \"5+93=98\", \"75+12=87\", \"15+47=62\"],\n  [\"97-85=12\", \"60-36=24\", \"1-0=1\", \"11+14=25\", \"2+34=36\"],\n  [\"39-5=34\", \"72-1=71\", \"30-27=3\", \"66-14=52\", \"26+37=63\"],\n  [\"92-18=74\", \"95-63=32\", \"46-10=36\", \"21+34=55\", \"21+32=53\"],\n  [\"67+2=69\", \"74-47=27\", \"99-56=43\", \"2+71=73\", \"22-20=2\"],\n  [\"32+64=96\", \"29-22=7\", \"87-51=36\", \"8+41=49\", \"93-36=57\"],\n  [\"37+10=47\", \"26+32=58\", \"33-20=13\", \"91-91=0\", \"79-29=50\"],\n  [\"95-30=65\", \"17+11=28\", \"82+4=86\", \"16+63=79\", \"52-48=4\"],\n  [\"12+70=82\", \"49+10=59\", \"17+16=33\", \"45+32=77\", \"92-77=15\"],\n  [\"65+14=79\", \"91+8=99\", \"56+1=57\", \"54+27=81\", \"11+53=64\"]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph of the document).\n$d.Paragraphs.Item(1).Range.Text = \"2025-09-29 Monday\"\n\n# Update each cell of the 20x5 addition/subtraction table in place,\n# preserving each cell paragraph/run formatting.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"69-43=26\"\n$t.Cell(1, 2).Range.Text = \"33+53=86\"\n$t.Cell(1, 3).Range.Text = \"17+53=70\"\n$t.Cell(1, 4).Range.Text = \"60+36=96\"\n$t.Cell(1, 5).Range.Text = \"18+32=50\"\n$t.Cell(2, 1).Range.Text = \"36+63=99\"\n$t.Cell(2, 2).Range.Text = \"64+29=93\"\n$t.Cell(2, 3).Range.Text = \"37-15=22\"\n$t.Cell(2, 4).Range.Text = \"97-32=65\"\n$t.Cell(2, 5).Range.Text = \"98-18=80\"\n$t.Cell(3, 1).Range.Text = \"26+3=29\"\n$t.Cell(3, 2).Range.Text = \"37-30=7\"\n$t.Cell(3, 3).Range.Text = \"5+88=93\"\n$t.Cell(3, 4).Range.Text = \"59-15=44\"\n$t.Cell(3, 5).Range.Text = \"26-15=11\"\n$t.Cell(4, 1).Range.Text = \"97+0=97\"\n$t.Cell(4, 2).Range.Text = \"3+93=96\"\n$t.Cell(4, 3).Range.Text = \"6+31=37\"\n$t.Cell(4, 4).Range.Text = \"77-49=28\"\n$t.Cell(4, 5).Range.Text = \"11+19=30\"\n$t.Cell(5, 1).Range.Text = \"15+18=33\"\n$t.Cell(5, 2).Range.Text = \"2-0=2\"\n$t.Cell(5, 3).Range.Text = \"53+33=86\"\n$t.Cell(5, 4).Range.Text = \"8+20=28\"\n$t.Cell(5, 5).Range.Text = \"76+5=81\"\n$t.Cell(6, 1).Range.Text = \"50+26=76\"\n$t.Cell(6, 2).Range.Text = \"57+2=59\"\n$t.Cell(6, 3).Range.Text = \"46+5=51\"\n$t.Cell(6, 4).Range.Text = \"33+46=79\"\n$t.Cell(6, 5).Range.Text = \"0+8=8\"\n$t.Cell(7, 1).Range.Text = \"94-82=12\"\n$t.Cell(7, 2).Range.Text = \"42-3=39\"\n$t.Cell(7, 3).Range.Text = \"58+19=77\"\n$t.Cell(7, 4).Range.Text = \"37-13=24\"\n$t.Cell(7, 5).Range.Text = \"19+80=99\"\n$t.Cell(8, 1).Range.Text = \"32-18=14\"\n$t.Cell(8, 2).Range.Text = \"46+5=51\"\n$t.Cell(8, 3).Range.Text = \"81-40=41\"\n$t.Cell(8, 4).Range.Text = \"30+67=97\"\n$t.Cell(8, 5).Range.Text = \"99-55=44\"\n$t.Cell(9, 1).Range.Text = \"69+21=90\"\n$t.Cell(9, 2).Range.Text = \"56+37=93\"\n$t.Cell(9, 3).Range.Text = \"49+22=71\"\n$t.Cell(9, 4).Range.Text = \"96-5=91\"\n$t.Cell(9, 5).Range.Text = \"37-8=29\"\n$t.Cell(10, 1).Range.Text = \"12+69=81\"\n$t.Cell(10, 2).Range.Text = \"6+2=8\"\n$t.Cell(10, 3).Range.Text = \"53+18=71\"\n$t.Cell(10, 4).Range.Text = \"35+16=51\"\n$t.Cell(10, 5).Range.Text = \"26+66=92\"\n$t.Cell(11, 1).Range.Text = \"59-43=16\"\n$t.Cell(11, 2).Range.Text = \"85-32=53\"\n$t.Cell(11, 3).Range.Text = \"5+93=98\"\n$t.Cell(11, 4).Range.Text = \"75+12=87\"\n$t.Cell(11, 5).Range.Text = \"15+47=62\"\n$t.Cell(12, 1).Range.Text = \"97-85=12\"\n$t.Cell(12, 2).Range.Text = \"60-36=24\"\n$t.Cell(12, 3).Range.Text = \"1-0=1\"\n$t.Cell(12, 4).Range.Text = \"11+14=25\"\n$t.Cell(12, 5).Range.Text = \"2+34=36\"\n$t.Cell(13, 1).Range.Text = \"39-5=34\"\n$t.Cell(13, 2).Range.Text = \"72-1=71\"\n$t.Cell(13, 3).Range.Text = \"30-27=3\"\n$t.Cell(13, 4).Range.Text = \"66-14=52\"\n$t.Cell(13, 5).Range.Text = \"26+37=63\"\n$t.Cell(14, 1).Range.Text = \"92-18=74\"\n$t.Cell(14, 2).Range.Text = \"95-63=32\"\n$t.Cell(14, 3).Range.Text = \"46-10=36\"\n$t.Cell(14, 4).Range.Text = \"21+34=55\"\n$t.Cell(14, 5).Range.Text = \"21+32=53\"\n$t.Cell(15, 1).Range.Text = \"67+2=69\"\n$t.Cell(15, 2).Range.Text = \"74-47=27\"\n$t.Cell(15, 3).Range.Text = \"99-56=43\"\n$t.Cell(15, 4).Range.Text = \"2+71=73\"\n$t.Cell(15, 5).Range.Text = \"22-20=2\"\n$t.Cell(16, 1).Range.Text = \"32+64=96\"\n$t.Cell(16, 2).Range.Text = \"29-22=7\"\n$t.Cell(16, 3).Range.Text = \"87-51=36\"\n$t.Cell(16, 4).Range.Text = \"8+41=49\"\n$t.Cell(16, 5).Range.Text = \"93-36=57\"\n$t.Cell(17, 1).Range.Text = \"37+10=47\"\n$t.Cell(17, 2).Range.Text = \"26+32=58\"\n$t.Cell(17, 3).Range.Text = \"33-20=13\"\n$t.Cell(17, 4).Range.Text = \"91-91=0\"\n$t.Cell(17, 5).Range.Text = \"79-29=50\"\n$t.Cell(18, 1).Range.Text = \"95-30=65\"\n$t.Cell(18, 2).Range.Text = \"17+11=28\"\n$t.Cell(18, 3).Range.Text = \"82+4=86\"\n$t.Cell(18, 4).Range.Text = \"16+63=79\"\n$t.Cell(18, 5).Range.Text = \"52-48=4\"\n$t.Cell(19, 1).Range.Text = \"12+70=82\"\n$t.Cell(19, 2).Range.Text = \"49+10=59\"\n$t.Cell(19, 3).Range.Text = \"17+16=33\"\n$t.Cell(19, 4).Range.Text = \"45+32=77\"\n$t.Cell(19, 5).Range.Text = \"92-77=15\"\n$t.Cell(20, 1).Range.Text = \"65+14=79\"\n$t.Cell(20, 2).Range.Text = \"91+8=99\"\n$t.Cell(20, 3).Range.Text = \"56+1=57\"\n$t.Cell(20, 4).Range.Text = \"54+27=81\"\n$t.Cell(20, 5).Range.Text = \"11+53=64\"\n"}
